$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.24610424041748
$ws.Range("B1").Value = 4.942692756652832
$ws.Range("C1").Value = 3.258103370666504
$ws.Range("D1").Value = 1.755290865898132
$ws.Range("E1").Value = 1.325338125228882
